# Update countries & provincias Spain
# - Senegal's case counts grew enough to overtake "Republica de Yibuti" and
#   "Georgia" in the ranked list (sheet is sorted descending by Casos
#   totales), so those three rows swap places: Senegal moves up to row 110,
#   Yibuti drops to row 111 and Georgia drops to row 112 (each keeping its
#   own previous numbers except Senegal, which gets fresh totals).
# - A handful of other countries got refreshed totals for this data pull.
# - The "updated at" timestamp footer cell moves from 14:52 to 15:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-rank Senegal / Republica de Yibuti / Georgia (rows 110-112) ---
$ws.Range("A110").Value = "Senegal"
$ws.Range("B110").Value = 299
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 183
$ws.Range("E110").Value = 114
$ws.Range("F110").Value = 1
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

$ws.Range("A111").Value = "Republica de Yibuti"
$ws.Range("B111").Value = 298
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 41
$ws.Range("E111").Value = 255
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 2

$ws.Range("A112").Value = "Georgia"
$ws.Range("B112").Value = 296
$ws.Range("C112").Value = 24
$ws.Range("D112").Value = 68
$ws.Range("E112").Value = 225
$ws.Range("F112").Value = 6
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 3

# --- Updated totals for other countries ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 587597
$ws.Range("C4").Value = 656
$ws.Range("D4").Value = 37202
$ws.Range("E4").Value = 526746
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 23649

# Alemania (row 8)
$ws.Range("B8").Value = 130383
$ws.Range("C8").Value = 311
$ws.Range("E8").Value = 58968
$ws.Range("G8").Value = 21
$ws.Range("H8").Value = 3215

# Brasil (row 17)
$ws.Range("B17").Value = 23830
$ws.Range("C17").Value = 400
$ws.Range("E17").Value = 19496

# Austria (row 20)
$ws.Range("B20").Value = 14159
$ws.Range("C20").Value = 118
$ws.Range("E20").Value = 6142

# Serbia (row 43)
$ws.Range("B43").Value = 4465
$ws.Range("C43").Value = 411
$ws.Range("E43").Value = 3971
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 94

# Finlandia (row 50)
$ws.Range("E50").Value = 2797
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 64

# Argentina (row 54)
$ws.Range("D54").Value = 559
$ws.Range("E54").Value = 1617
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 101

# Kazajistan (row 70)
$ws.Range("D70").Value = 203
$ws.Range("E70").Value = 985

# Azerbaiyan (row 71)
$ws.Range("B71").Value = 1197
$ws.Range("C71").Value = 49
$ws.Range("D71").Value = 351
$ws.Range("E71").Value = 833
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 13

# Camerun (row 79)
$ws.Range("D79").Value = 130
$ws.Range("E79").Value = 704
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 14

# Sri Lanka (row 117)
$ws.Range("D117").Value = 61
$ws.Range("E117").Value = 151

# Togo (row 139)
$ws.Range("D139").Value = 32
$ws.Range("E139").Value = 42

# Laos (row 173)
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 18

# Republica de Africa Central (row 194)
$ws.Range("D194").Value = 4
$ws.Range("E194").Value = 7

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 15:22"
